$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Title (B) and Year (D) values for rows 2-10, leaving the
# bullet markers in column A untouched. This mirrors the already-empty
# rows further down the sheet (e.g. row 11+) where B has no text and D
# has no cell at all.
$ws.Range("B2:B10").ClearContents()
$ws.Range("D2:D10").ClearContents()
